$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 779480
$ws.Range("I31").Value = 779480
$ws.Range("K31").Value = 2338440
$ws.Range("M31").Value = -2338210
$ws.Range("H41").Value = 211.08333
$ws.Range("I41").Value = 205.83333
$ws.Range("J41").Value = 216.33333
$ws.Range("K41").Value = 205.83333
$ws.Range("L41").Value = 216.33333
$ws.Range("M41").Value = 234.16667
$ws.Range("N41").Value = -1096.33333
$ws.Range("H101").Value = 1208.3334
$ws.Range("I101").Value = 330
$ws.Range("J101").Value = 2965
$ws.Range("K101").Value = 990
$ws.Range("L101").Value = 8895
$ws.Range("M101").Value = 632
$ws.Range("N101").Value = -12139
$ws.Range("H112").Value = 2188.9614
$ws.Range("J112").Value = 2357.087
$ws.Range("L112").Value = 7071.261
$ws.Range("N112").Value = -9287.261
$ws.Range("H116").Value = 16669620
$ws.Range("I116").Value = 40002770
$ws.Range("J116").Value = 3085.1428
$ws.Range("K116").Value = 40002770
$ws.Range("L116").Value = 3085.1428
$ws.Range("M116").Value = -39999328
$ws.Range("N116").Value = -9969.1428

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 37970.57
$ws.Range("I86").Value = 2412.9375
$ws.Range("J86").Value = 85380.75
$ws.Range("K86").Value = 2412.9375
$ws.Range("L86").Value = 85380.75
$ws.Range("M86").Value = -1289.9375
$ws.Range("N86").Value = -87626.75
$ws.Range("H89").Value = 37970.57
$ws.Range("I89").Value = 2412.9375
$ws.Range("J89").Value = 85380.75
$ws.Range("K89").Value = 12064.6875
$ws.Range("L89").Value = 426903.75
$ws.Range("M89").Value = -6448.6875
$ws.Range("N89").Value = -438135.75
$ws.Range("H108").Value = 36499.668
$ws.Range("J108").Value = 36499.668
$ws.Range("L108").Value = 36499.668
$ws.Range("N108").Value = -44179.668

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1933.5769
$ws.Range("I58").Value = 1584.3125
$ws.Range("J58").Value = 2492.4
$ws.Range("K58").Value = 1584.3125
$ws.Range("L58").Value = 2492.4
$ws.Range("M58").Value = -1381.3125
$ws.Range("N58").Value = -2898.4
$ws.Range("H99").Value = 1920.5714
$ws.Range("I99").Value = 2254.7778
$ws.Range("J99").Value = 1319
$ws.Range("K99").Value = 2254.7778
$ws.Range("L99").Value = 1319
$ws.Range("M99").Value = -756.7777999999998
$ws.Range("N99").Value = -4315
$ws.Range("H114").Value = 43683.668
$ws.Range("J114").Value = 43683.668
$ws.Range("L114").Value = 43683.668
$ws.Range("N114").Value = -52361.668
$ws.Range("H126").Value = 1920.5714
$ws.Range("I126").Value = 2254.7778
$ws.Range("J126").Value = 1319
$ws.Range("K126").Value = 6764.3334
$ws.Range("L126").Value = 3957
$ws.Range("M126").Value = -4294.3334
$ws.Range("N126").Value = -8897
$ws.Range("H136").Value = 1933.5769
$ws.Range("I136").Value = 1584.3125
$ws.Range("J136").Value = 2492.4
$ws.Range("K136").Value = 4752.9375
$ws.Range("L136").Value = 7477.200000000001
$ws.Range("M136").Value = -2202.9375
$ws.Range("N136").Value = -12577.2

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1787.0869
$ws.Range("J5").Value = 999.95
$ws.Range("L5").Value = 2999.85
$ws.Range("N5").Value = -3223.85
$ws.Range("H122").Value = 1004.5
$ws.Range("J122").Value = 1555
$ws.Range("L122").Value = 13995
$ws.Range("N122").Value = -18895
$ws.Range("H135").Value = 1787.0869
$ws.Range("J135").Value = 999.95
$ws.Range("L135").Value = 8999.550000000001
$ws.Range("N135").Value = -14069.55

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7829.8335
$ws.Range("I70").Value = 8719.75
$ws.Range("J70").Value = 6050
$ws.Range("K70").Value = 8719.75
$ws.Range("L70").Value = 6050
$ws.Range("M70").Value = -8449.75
$ws.Range("N70").Value = -6590
$ws.Range("H73").Value = 7829.8335
$ws.Range("I73").Value = 8719.75
$ws.Range("J73").Value = 6050
$ws.Range("K73").Value = 8719.75
$ws.Range("L73").Value = 6050
$ws.Range("M73").Value = -7783.75
$ws.Range("N73").Value = -7922
$ws.Range("H108").Value = 39000
$ws.Range("J108").Value = 39000
$ws.Range("L108").Value = 39000
$ws.Range("N108").Value = -46680
$ws.Range("H113").Value = 1826.9412
$ws.Range("I113").Value = 1235.2307
$ws.Range("J113").Value = 3750
$ws.Range("K113").Value = 1235.2307
$ws.Range("L113").Value = 3750
$ws.Range("M113").Value = 934.7692999999999
$ws.Range("N113").Value = -8090

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 12889.263
$ws.Range("I61").Value = 15686.4
$ws.Range("J61").Value = 2400
$ws.Range("K61").Value = 15686.4
$ws.Range("L61").Value = 2400
$ws.Range("M61").Value = -15484.4
$ws.Range("N61").Value = -2804
$ws.Range("H100").Value = 15495
$ws.Range("I100").Value = 26990
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 26990
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -26449
$ws.Range("N100").Value = -5082
$ws.Range("H113").Value = 12889.263
$ws.Range("I113").Value = 15686.4
$ws.Range("J113").Value = 2400
$ws.Range("K113").Value = 15686.4
$ws.Range("L113").Value = 2400
$ws.Range("M113").Value = -13516.4
$ws.Range("N113").Value = -6740

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 56000
$ws.Range("H73").Value = 56000
$ws.Range("H96").Value = 2315.6
$ws.Range("I96").Value = 2000
$ws.Range("J96").Value = 2526
$ws.Range("K96").Value = 2000
$ws.Range("L96").Value = 2526
$ws.Range("M96").Value = -627
$ws.Range("N96").Value = -5272
$ws.Range("H100").Value = 19018
$ws.Range("I100").Value = 33783
$ws.Range("J100").Value = 1300
$ws.Range("K100").Value = 67566
$ws.Range("L100").Value = 2600
$ws.Range("M100").Value = -67025
$ws.Range("N100").Value = -3682
$ws.Range("H136").Value = 3369.5881
$ws.Range("I136").Value = 2754.1428
$ws.Range("J136").Value = 3800.4
$ws.Range("K136").Value = 8262.428400000001
$ws.Range("L136").Value = 11401.2
$ws.Range("M136").Value = -5712.428400000001
$ws.Range("N136").Value = -16501.2
